$wb = $excel.ActiveWorkbook

# The workbook contains duplicated data across the "展览" and "全部类型" sheets.
# Update the "想去人数" (F column) counts on both sheets to match the new values.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 132
    $ws.Range("F9").Value = 593
    $ws.Range("F10").Value = 407
}
